$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 857.5
$ws.Range("I9").Value = 1001
$ws.Range("J9").Value = 570.5
$ws.Range("K9").Value = 1001
$ws.Range("L9").Value = 570.5
$ws.Range("M9").Value = -832
$ws.Range("N9").Value = -908.5
$ws.Range("H33").Value = 102.454544
$ws.Range("I33").Value = 120.875
$ws.Range("K33").Value = 120.875
$ws.Range("M33").Value = 108.125
$ws.Range("H59").Value = 4999.6
$ws.Range("J59").Value = 6666.3335
$ws.Range("L59").Value = 19999.0005
$ws.Range("N59").Value = -21113.0005
$ws.Range("H88").Value = 2562.25
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2562.25
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2562.25
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -3374.25
$ws.Range("H91").Value = 2562.25
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2562.25
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2562.25
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -5370.25
$ws.Range("H132").Value = 2268.3
$ws.Range("I132").Value = 2226
$ws.Range("J132").Value = 2649
$ws.Range("K132").Value = 6678
$ws.Range("L132").Value = 7947
$ws.Range("M132").Value = -4148
$ws.Range("N132").Value = -13007
$ws.Range("H135").Value = 419.84616
$ws.Range("I135").Value = 486.9
$ws.Range("K135").Value = 4382.099999999999
$ws.Range("M135").Value = -1847.099999999999
$ws.Range("H138").Value = 4466
$ws.Range("J138").Value = 4775.7334
$ws.Range("L138").Value = 14327.2002
$ws.Range("N138").Value = -24607.2002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5744.5415
$ws.Range("I32").Value = 2724.7144
$ws.Range("J32").Value = 26883.334
$ws.Range("K32").Value = 2724.7144
$ws.Range("L32").Value = 26883.334
$ws.Range("M32").Value = -2437.7144
$ws.Range("N32").Value = -27457.334
$ws.Range("H36").Value = 5048.875
$ws.Range("I36").Value = 6316
$ws.Range("J36").Value = 1247.5
$ws.Range("K36").Value = 6316
$ws.Range("L36").Value = 1247.5
$ws.Range("M36").Value = -5970
$ws.Range("N36").Value = -1939.5
$ws.Range("H45").Value = 4246
$ws.Range("I45").Value = 1444
$ws.Range("J45").Value = 4868.6665
$ws.Range("K45").Value = 1444
$ws.Range("L45").Value = 4868.6665
$ws.Range("M45").Value = -1067
$ws.Range("N45").Value = -5622.6665
$ws.Range("H61").Value = 3124.5
$ws.Range("I61").Value = 2250
$ws.Range("J61").Value = 3999
$ws.Range("K61").Value = 2250
$ws.Range("L61").Value = 3999
$ws.Range("M61").Value = -2038
$ws.Range("N61").Value = -4423
$ws.Range("H97").Value = 858.4666999999999
$ws.Range("I97").Value = 406.41666
$ws.Range("K97").Value = 406.41666
$ws.Range("M97").Value = 89.58334000000002
$ws.Range("H136").Value = 3124.5
$ws.Range("I136").Value = 2250
$ws.Range("J136").Value = 3999
$ws.Range("K136").Value = 6750
$ws.Range("L136").Value = 11997
$ws.Range("M136").Value = -4200
$ws.Range("N136").Value = -17097

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2281.8667
$ws.Range("I20").Value = 1553.2
$ws.Range("J20").Value = 3739.2
$ws.Range("K20").Value = 1553.2
$ws.Range("L20").Value = 3739.2
$ws.Range("M20").Value = -1306.2
$ws.Range("N20").Value = -4233.2
$ws.Range("H22").Value = 437.5
$ws.Range("J22").Value = 290
$ws.Range("L22").Value = 290
$ws.Range("N22").Value = -636
$ws.Range("H29").Value = 3016
$ws.Range("I29").Value = 3016
$ws.Range("K29").Value = 3016
$ws.Range("M29").Value = -2727
$ws.Range("H82").Value = 29102.4
$ws.Range("J82").Value = 60000
$ws.Range("L82").Value = 60000
$ws.Range("N82").Value = -60766
$ws.Range("H85").Value = 29102.4
$ws.Range("J85").Value = 60000
$ws.Range("L85").Value = 60000
$ws.Range("N85").Value = -62652
$ws.Range("H99").Value = 2458.1667
$ws.Range("I99").Value = 2709.8
$ws.Range("K99").Value = 2709.8
$ws.Range("M99").Value = -1211.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 571621.1
$ws.Range("I6").Value = 571621.1
$ws.Range("K6").Value = 571621.1
$ws.Range("M6").Value = -571508.1
$ws.Range("H22").Value = 1675
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H31").Value = 5750
$ws.Range("I31").Value = 5000
$ws.Range("K31").Value = 5000
$ws.Range("M31").Value = -4705
$ws.Range("H34").Value = 5750
$ws.Range("I34").Value = 5000
$ws.Range("K34").Value = 5000
$ws.Range("M34").Value = -4798
$ws.Range("H105").Value = 1472.625
$ws.Range("I105").Value = 1056.6
$ws.Range("J105").Value = 2166
$ws.Range("K105").Value = 1056.6
$ws.Range("L105").Value = 2166
$ws.Range("M105").Value = 690.4000000000001
$ws.Range("N105").Value = -5660

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 777.36365
$ws.Range("I34").Value = 594.1667
$ws.Range("J34").Value = 997.2
$ws.Range("K34").Value = 1782.5001
$ws.Range("L34").Value = 2991.6
$ws.Range("M34").Value = -1698.5001
$ws.Range("N34").Value = -3159.6
$ws.Range("H55").Value = 788.6667
$ws.Range("J55").Value = 1400
$ws.Range("L55").Value = 4200
$ws.Range("N55").Value = -4554
$ws.Range("H131").Value = 857.5909
$ws.Range("J131").Value = 990
$ws.Range("L131").Value = 2970
$ws.Range("N131").Value = -13050

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H24").Value = 4000000
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H136").Value = 20956.5
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 20956.5
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 62869.5
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -67969.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1651.1818
$ws.Range("I16").Value = 1651.1818
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1651.1818
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1481.1818
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 1542.9286
$ws.Range("I22").Value = 1509.2727
$ws.Range("J22").Value = 1666.3334
$ws.Range("K22").Value = 1509.2727
$ws.Range("L22").Value = 1666.3334
$ws.Range("M22").Value = -1214.2727
$ws.Range("N22").Value = -2256.3334
$ws.Range("H27").Value = 1542.9286
$ws.Range("I27").Value = 1509.2727
$ws.Range("J27").Value = 1666.3334
$ws.Range("K27").Value = 1509.2727
$ws.Range("L27").Value = 1666.3334
$ws.Range("M27").Value = -1402.2727
$ws.Range("N27").Value = -1880.3334
$ws.Range("H29").Value = 50000000
$ws.Range("I29").Value = 50000000
$ws.Range("K29").Value = 50000000
$ws.Range("M29").Value = -49999705
$ws.Range("H68").Value = 2178.9
$ws.Range("I68").Value = 1639
$ws.Range("J68").Value = 2718.8
$ws.Range("K68").Value = 1639
$ws.Range("L68").Value = 2718.8
$ws.Range("M68").Value = -890
$ws.Range("N68").Value = -4216.8
$ws.Range("H71").Value = 2178.9
$ws.Range("I71").Value = 1639
$ws.Range("J71").Value = 2718.8
$ws.Range("K71").Value = 8195
$ws.Range("L71").Value = 13594
$ws.Range("M71").Value = -4451
$ws.Range("N71").Value = -21082
$ws.Range("H132").Value = 4525.4614
$ws.Range("I132").Value = 4703.5557
$ws.Range("J132").Value = 4124.75
$ws.Range("K132").Value = 14110.6671
$ws.Range("L132").Value = 12374.25
$ws.Range("M132").Value = -11580.6671
$ws.Range("N132").Value = -17434.25
$ws.Range("H136").Value = 3223.5715
$ws.Range("J136").Value = 4200
$ws.Range("L136").Value = 12600
$ws.Range("N136").Value = -17700

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 12500
$ws.Range("J81").Value = 12500
$ws.Range("L81").Value = 25000
$ws.Range("N81").Value = -27122
$ws.Range("H84").Value = 12500
$ws.Range("J84").Value = 12500
$ws.Range("L84").Value = 125000
$ws.Range("N84").Value = -135608
$ws.Range("H126").Value = 1749.5
$ws.Range("I126").Value = 1749.5
$ws.Range("K126").Value = 5248.5
$ws.Range("M126").Value = -2778.5
